$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 is being updated with a new patient record. B13/C13/D13 reuse text
# that already exists verbatim on row 12 (LUIGGI / PASACHE LOPERA / 21), so a
# straight value copy (Copy + PasteSpecial values) carries over the original
# shared-string cell type without Excel re-inferring a numeric type for the
# digit-only strings.
$ws.Range("B12:D12").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

# A13's new code ("11493875") is brand new text and doesn't exist anywhere in
# the sheet yet to copy from. Building it via a formula that concatenates
# text literals forces Excel to treat the result as a STRING (never as a
# number), so copying that formula's cached value across keeps it text too -
# unlike assigning the digits straight to .Value, which Excel would silently
# reinterpret as a numeric cell.
$scratch = $ws.Range("H1")
$scratch.Formula = '=""&"11493875"'
$scratch.Copy()
$ws.Range("A13").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false
$scratch.ClearContents()
